# Update latest output (run 112)
# - Schedule sheet: refresh computed Cost / Unit Cost for rows 4 and 5
# - Detailed sheet: refresh Price forecasts/historicals for rows 45-96
#   (includes a few rows flipping Type from "forecast" to "historical")

$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet updates ---
$wsSchedule.Range("E4").Value2 = 448.0351589999999
$wsSchedule.Range("F4").Value2 = 29.63195496031745
$wsSchedule.Range("E5").Value2 = -35.15649150000001
$wsSchedule.Range("F5").Value2 = -1.033406569664903

# --- Detailed sheet updates ---
$wsDetailed.Range("B45").Value2 = 65
$wsDetailed.Range("B46").Value2 = 61.69483
$wsDetailed.Range("B47").Value2 = 58.06057
$wsDetailed.Range("C47").Value = "historical"
$wsDetailed.Range("B48").Value2 = 57.98309
$wsDetailed.Range("C48").Value = "historical"
$wsDetailed.Range("B49").Value2 = 59.16
$wsDetailed.Range("C49").Value = "historical"
$wsDetailed.Range("B51").Value2 = 57.08
$wsDetailed.Range("B59").Value2 = 60.66677
$wsDetailed.Range("B60").Value2 = 64.11489
$wsDetailed.Range("B61").Value2 = 73.92341
$wsDetailed.Range("B62").Value2 = 77.09846
$wsDetailed.Range("B63").Value2 = 66.23294
$wsDetailed.Range("B64").Value2 = 55.87148
$wsDetailed.Range("B65").Value2 = 8.63298
$wsDetailed.Range("B66").Value2 = 1.71524
$wsDetailed.Range("B67").Value2 = 0.7
$wsDetailed.Range("B68").Value2 = -0.79345
$wsDetailed.Range("B69").Value2 = -5.58973
$wsDetailed.Range("B70").Value2 = -7.00424
$wsDetailed.Range("B71").Value2 = -6.67637
$wsDetailed.Range("B72").Value2 = -7.88167
$wsDetailed.Range("B73").Value2 = -7.96339
$wsDetailed.Range("B74").Value2 = -7.78102
$wsDetailed.Range("B75").Value2 = -8.56077
$wsDetailed.Range("B76").Value2 = -9.99
$wsDetailed.Range("B77").Value2 = -8.33366
$wsDetailed.Range("B78").Value2 = -9.494960000000001
$wsDetailed.Range("B79").Value2 = -8.280570000000001
$wsDetailed.Range("B80").Value2 = -7.82781
$wsDetailed.Range("B83").Value2 = -5.86515
$wsDetailed.Range("B84").Value2 = -4.69749
$wsDetailed.Range("B85").Value2 = 5.15736
$wsDetailed.Range("B86").Value2 = 9.813750000000001
$wsDetailed.Range("B87").Value2 = 33.01451
$wsDetailed.Range("B88").Value2 = 56.98
$wsDetailed.Range("B90").Value2 = 57.3
$wsDetailed.Range("B91").Value2 = 57.3
$wsDetailed.Range("B93").Value2 = 57.09607
$wsDetailed.Range("B94").Value2 = 56.98
$wsDetailed.Range("B95").Value2 = 57.06007
$wsDetailed.Range("B96").Value2 = 57.06007
